$wb = $excel.ActiveWorkbook

# --- Insert a new worksheet "StuffToImport" before "Models to make" -------
# (sheet references in this bridge are positional, so grab a fresh handle
#  to "Models to make" again *after* the insertion shifts indices)
$newSheet = $wb.Worksheets.Add($wb.Worksheets.Item("Models to make"))
$newSheet.Name = "StuffToImport"
$modelsSheet = $wb.Worksheets.Item("Models to make")

# --- BaseTable (A1:C4) ------------------------------------------------
$newSheet.Range("A1").Value = "Base"
$newSheet.Range("B1").Value = "Durability"
$newSheet.Range("C1").Value = "Speed"

$newSheet.Range("A2").Value = "WoodenBase"
$newSheet.Range("B2").Value = 10
$newSheet.Range("C2").Value = 5

$newSheet.Range("A3").Value = "MetalBase"
$newSheet.Range("B3").Value = 75
$newSheet.Range("C3").Value = 7

$newSheet.Range("A4").Value = "LeatherBase"
$newSheet.Range("B4").Value = 25
$newSheet.Range("C4").Value = 15

# --- TopTable (E1:G4) ---------------------------------------------------
$newSheet.Range("E1").Value = "Top"
$newSheet.Range("F1").Value = "Durability"
$newSheet.Range("G1").Value = "Speed"

$newSheet.Range("E2").Value = "WoodenHook"
$newSheet.Range("F2").Value = 5
$newSheet.Range("G2").Value = 5

$newSheet.Range("E3").Value = "MetalHook"
$newSheet.Range("F3").Value = 10
$newSheet.Range("G3").Value = 7

$newSheet.Range("E4").Value = "MetalHand"
$newSheet.Range("F4").Value = 15
$newSheet.Range("G4").Value = 8

# --- Name / LeatherBase_MetalHand note ----------------------------------
$newSheet.Range("A7").Value = "Name"
$newSheet.Range("A8").Value = "LeatherBase_MetalHand"

# --- Turn the two ranges into real tables (ListObjects) -----------------
$baseTable = $newSheet.ListObjects.Add(1, $newSheet.Range("A1:C4"), $null, 1)
$baseTable.Name = "BaseTable"
$baseTable.TableStyle = "TableStyleLight1"

$topTable = $newSheet.ListObjects.Add(1, $newSheet.Range("E1:G4"), $null, 1)
$topTable.Name = "TopTable"
$topTable.TableStyle = "TableStyleLight1"

# --- Cosmetics on the new sheet -----------------------------------------
$newSheet.Columns.Item(1).ColumnWidth = 11.417
$newSheet.Columns.Item(2).ColumnWidth = 8.417
$newSheet.Columns.Item(5).ColumnWidth = 10.584
$newSheet.Columns.Item(6).ColumnWidth = 8.584
$newSheet.Columns.Item(7).ColumnWidth = 11.75
$newSheet.PageSetup.Orientation = 1
$newSheet.Activate()
$newSheet.Range("E1").Select()
$excel.ActiveWindow.Zoom = 205

# --- "Models to make" sheet: widen column A ------------------------------
$modelsSheet.Columns.Item(1).ColumnWidth = 22.917
